# Apply scheduled Typhon_Profits valuation updates across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5101658.5
$ws.Range("J17").Value = 5370067
$ws.Range("L17").Value = 16110201
$ws.Range("N17").Value = -16110537

$ws.Range("H62").Value = 8258.817999999999
$ws.Range("I62").Value = 6491
$ws.Range("J62").Value = 9732
$ws.Range("K62").Value = 6491
$ws.Range("L62").Value = 9732
$ws.Range("M62").Value = -5867
$ws.Range("N62").Value = -10980

$ws.Range("H65").Value = 8258.817999999999
$ws.Range("I65").Value = 6491
$ws.Range("J65").Value = 9732
$ws.Range("K65").Value = 32455
$ws.Range("L65").Value = 48660
$ws.Range("M65").Value = -29335
$ws.Range("N65").Value = -54900

$ws.Range("H100").Value = 125002800
$ws.Range("I100").Value = 250001870
$ws.Range("J100").Value = 3725
$ws.Range("K100").Value = 250001870
$ws.Range("L100").Value = 3725
$ws.Range("M100").Value = -250001329
$ws.Range("N100").Value = -4807

$ws.Range("H121").Value = 12359.5
$ws.Range("J121").Value = 15249.375
$ws.Range("L121").Value = 45748.125
$ws.Range("N121").Value = -49242.125

$ws.Range("H129").Value = 233512.9
$ws.Range("J129").Value = 295266.97
$ws.Range("L129").Value = 885800.9099999999
$ws.Range("N129").Value = -895800.9099999999

$ws.Range("H132").Value = 2160.3845
$ws.Range("I132").Value = 2330.182
$ws.Range("J132").Value = 1226.5
$ws.Range("K132").Value = 6990.545999999999
$ws.Range("L132").Value = 3679.5
$ws.Range("M132").Value = -4460.545999999999
$ws.Range("N132").Value = -8739.5

$ws.Range("H141").Value = 3647.875
$ws.Range("I141").Value = 2836.6
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 8509.799999999999
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -3329.799999999999
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 66668976
$ws.Range("I74").Value = 142857810
$ws.Range("K74").Value = 142857810
$ws.Range("M74").Value = -142856936

$ws.Range("H77").Value = 66668976
$ws.Range("I77").Value = 142857810
$ws.Range("K77").Value = 714289050
$ws.Range("M77").Value = -714284682

$ws.Range("H102").Value = 1567.1428
$ws.Range("I102").Value = 1514
$ws.Range("J102").Value = 1700
$ws.Range("K102").Value = 1514
$ws.Range("L102").Value = 1700
$ws.Range("M102").Value = 108
$ws.Range("N102").Value = -4944

$ws.Range("H122").Value = 2161.6667
$ws.Range("I122").Value = 1822.0555
$ws.Range("J122").Value = 4199.3335
$ws.Range("K122").Value = 5466.166499999999
$ws.Range("L122").Value = 12598.0005
$ws.Range("M122").Value = -3016.166499999999
$ws.Range("N122").Value = -17498.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 655.7778
$ws.Range("I94").Value = 457.42856
$ws.Range("J94").Value = 1350
$ws.Range("K94").Value = 457.42856
$ws.Range("L94").Value = 1350
$ws.Range("M94").Value = -6.428560000000004
$ws.Range("N94").Value = -2252

$ws.Range("H132").Value = 40000
$ws.Range("J132").Value = 40000
$ws.Range("L132").Value = 40000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 55559050
$ws.Range("I62").Value = 125003890
$ws.Range("J62").Value = 3174.8
$ws.Range("K62").Value = 125003890
$ws.Range("L62").Value = 3174.8
$ws.Range("M62").Value = -125003266
$ws.Range("N62").Value = -4422.8

$ws.Range("H65").Value = 55559050
$ws.Range("I65").Value = 125003890
$ws.Range("J65").Value = 3174.8
$ws.Range("K65").Value = 625019450
$ws.Range("L65").Value = 15874
$ws.Range("M65").Value = -625016330
$ws.Range("N65").Value = -22114

$ws.Range("H132").Value = 2901.182
$ws.Range("I132").Value = 2069.3333
$ws.Range("J132").Value = 6644.5
$ws.Range("K132").Value = 6207.999899999999
$ws.Range("L132").Value = 19933.5
$ws.Range("M132").Value = -3677.999899999999
$ws.Range("N132").Value = -24993.5

$ws.Range("H134").Value = 1315.7693
$ws.Range("I134").Value = 1218.6364
$ws.Range("J134").Value = 1850
$ws.Range("K134").Value = 3655.9092
$ws.Range("L134").Value = 5550
$ws.Range("M134").Value = -1120.9092
$ws.Range("N134").Value = -10620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2037.7142
$ws.Range("J5").Value = 2652.5
$ws.Range("L5").Value = 7957.5
$ws.Range("N5").Value = -8181.5

$ws.Range("H121").Value = 1069.5217
$ws.Range("I121").Value = 299.5
$ws.Range("K121").Value = 898.5
$ws.Range("M121").Value = 411.5

$ws.Range("H131").Value = 703.9394
$ws.Range("J131").Value = 703.9394
$ws.Range("L131").Value = 2111.8182
$ws.Range("N131").Value = -12191.8182

$ws.Range("H135").Value = 2037.7142
$ws.Range("J135").Value = 2652.5
$ws.Range("L135").Value = 23872.5
$ws.Range("N135").Value = -28942.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6581.931
$ws.Range("I70").Value = 3082.913
$ws.Range("K70").Value = 3082.913
$ws.Range("M70").Value = -2812.913

$ws.Range("H73").Value = 6581.931
$ws.Range("I73").Value = 3082.913
$ws.Range("K73").Value = 3082.913
$ws.Range("M73").Value = -2146.913

$ws.Range("H132").Value = 57854.668
$ws.Range("I132").Value = 1256
$ws.Range("K132").Value = 3768
$ws.Range("M132").Value = -1238

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I16").Value = 424.65
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 424.65
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -254.65
$ws.Range("N16").ClearContents()

$ws.Range("H40").Value = 4091.2
$ws.Range("I40").Value = 2823.111
$ws.Range("J40").Value = 7352
$ws.Range("K40").Value = 2823.111
$ws.Range("L40").Value = 7352
$ws.Range("M40").Value = -2687.111
$ws.Range("N40").Value = -7624

$ws.Range("H136").Value = 1407.8572
$ws.Range("I136").Value = 1407.8572
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4223.571599999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1673.571599999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 52273240
$ws.Range("I107").Value = 90909240
$ws.Range("K107").Value = 272727720
$ws.Range("M107").Value = -272725800
